$d = $word.ActiveDocument

$pairs = @(
    @("85×22=1870", "74×40=2960"),
    @("50×67=3350", "64×34=2176"),
    @("52×18=936",  "81×85=6885"),
    @("72×59=4248", "67×12=804"),
    @("21×97=2037", "21×29=609"),
    @("75×67=5025", "68×63=4284"),
    @("22×29=638",  "13×47=611"),
    @("33×43=1419", "29×68=1972"),
    @("29×35=1015", "35×13=455"),
    @("47×86=4042", "71×73=5183"),
    @("75×15=1125", "16×83=1328"),
    @("18×60=1080", "28×41=1148"),
    @("54×55=2970", "49×49=2401"),
    @("42×34=1428", "61×95=5795"),
    @("93×13=1209", "30×71=2130"),
    @("27×90=2430", "34×88=2992"),
    @("49×79=3871", "51×45=2295"),
    @("57×72=4104", "59×37=2183"),
    @("47×17=799",  "95×21=1995"),
    @("45×68=3060", "37×27=999"),
    @("92×95=8740", "91×15=1365"),
    @("86×74=6364", "78×14=1092"),
    @("32×47=1504", "17×57=969"),
    @("69×63=4347", "91×96=8736"),
    @("69×29=2001", "29×21=609")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
